# Insert a new data row before row 62 (pushes existing rows 62:93 down to 63:94)
# and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(62).Insert()

$ws.Range("A62").Value = 10
$ws.Range("B62").Value = "Vega Modelo de Temuco"
$ws.Range("C62").Value = "La Araucanía"
$ws.Range("D62").Value = 44741
$ws.Range("E62").Value = 9
$ws.Range("F62").Value = 100114002
$ws.Range("G62").Value = "Camote"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 30
$ws.Range("K62").Value = 20000
$ws.Range("L62").Value = 20000
$ws.Range("M62").Value = 20000
$ws.Range("N62").Value = "$/malla 20 kilos"
$ws.Range("O62").Value = "Perú"
$ws.Range("P62").Value = 1000
$ws.Range("Q62").Value = 20
$ws.Range("R62").Value = "Hortaliza"
